$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add RoomNo-style distinguishing suffixes to the lab Course Code values
# (column A) so each lab section's code reflects its own room/section
# instead of sharing a duplicated code across A/B batches.
$ws.Range("A6").Value  = "CS102A-L1"
$ws.Range("A7").Value  = "CS102A-L2"
$ws.Range("A8").Value  = "CS102B-L1"
$ws.Range("A9").Value  = "CS102B-L2"
$ws.Range("A12").Value = "CS101A-L1"
$ws.Range("A13").Value = "CS101A-L2"
$ws.Range("A14").Value = "CS101B-L"
$ws.Range("A19").Value = "CS103A-L"
$ws.Range("A20").Value = "CS103B-L"
$ws.Range("A25").Value = "NS104A-L"
$ws.Range("A26").Value = "NS104B-L"
$ws.Range("A29").Value = "CS104A-L"
$ws.Range("A30").Value = "CS104B-L"

# Reset the view back to the top of the sheet (the saved workbook was
# scrolled to A44 with D59 selected beforehand).
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A1").Select()
